$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same
#    column layout/formatting), and place it right before "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newQ = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ.Name = "2022-Q1"

# Force the D:G columns to be stored as text (matching the rest of the
# workbook, which keeps these numeric-looking figures as plain text)
# before writing the new values, otherwise Excel will auto-convert the
# strings to numbers.
$newQ.Range("D2:G5").NumberFormat = "@"

$newQ.Cells.Item(2,4).Value = "0.57"
$newQ.Cells.Item(2,5).Value = "92.18"
$newQ.Cells.Item(2,6).Value = "9.06"
$newQ.Cells.Item(2,7).Value = "0.0516"
$newQ.Cells.Item(2,8).Value = 3

$newQ.Cells.Item(3,4).Value = "0.22"
$newQ.Cells.Item(3,5).Value = "91.35"
$newQ.Cells.Item(3,6).Value = "8.48"
$newQ.Cells.Item(3,7).Value = "0.0187"
$newQ.Cells.Item(3,8).Value = 2

$newQ.Cells.Item(4,4).Value = "0.05"
$newQ.Cells.Item(4,5).Value = "91.35"
$newQ.Cells.Item(4,6).Value = "8.48"
$newQ.Cells.Item(4,7).Value = "0.0042"
$newQ.Cells.Item(4,8).Value = 2

$newQ.Cells.Item(5,4).Value = "0.04"
$newQ.Cells.Item(5,5).Value = "92.18"
$newQ.Cells.Item(5,6).Value = "9.06"
$newQ.Cells.Item(5,7).Value = "0.0036"
$newQ.Cells.Item(5,8).Value = 3

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row for 2022-Q1 at
#    the top of the data (row 2) and shift the remaining rows down,
#    renumbering the index column (A) accordingly.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.08

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5

Write-Output "done"
